$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.1355019156032807
$ws.Range("D2").Value = 0.8934472299012921

# Row 3
$ws.Range("C3").Value = -0.2331079371900004
$ws.Range("D3").Value = 0.8178327860174979

# Row 4
$ws.Range("C4").Value = 0.3315013205270076
$ws.Range("D4").Value = 0.7434029926237158

# Row 5
$ws.Range("C5").Value = 2.130116879079089
$ws.Range("D5").Value = 0.04459022723629058
$ws.Range("G5").Value = "Sí"

# Row 6
$ws.Range("C6").Value = -0.4026746049520101
$ws.Range("D6").Value = 0.6910737157725855

# Row 7
$ws.Range("C7").Value = 0.1282869585720923
$ws.Range("D7").Value = 0.8990873935755945

# Row 8
$ws.Range("C8").Value = 1.574479472040786
$ws.Range("D8").Value = 0.1296488090973209

# Row 9
$ws.Range("C9").Value = 0.7368350409708748
$ws.Range("D9").Value = 0.4690082112456286

# Row 10
$ws.Range("C10").Value = 1.324868432970911
$ws.Range("D10").Value = 0.1988121911894549

# Row 11
$ws.Range("C11").Value = 1.342920992703276
$ws.Range("D11").Value = 0.1929894545432891
